$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-22 13:53:11"
$wsZhCn.Range("G3").Value = "2016-02-22 13:54:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-22 13:53:26"
$wsDeDe.Range("G3").Value = "2016-02-22 13:54:38"
